# Update res_bus/vm_pu.xlsx values for "case with 380 kV done"
# Slack bus voltage setpoint (column B) moved from 1.05 p.u. to 1.02 p.u.,
# and the resulting bus voltage magnitudes (columns C-F, I-N) were recalculated
# for every data row (Excel rows 2-25).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.020260089553636
$ws.Range("D2").Value = 1.031669768852508
$ws.Range("E2").Value = 1.021293454059137
$ws.Range("F2").Value = 1.032856848201286
$ws.Range("I2").Value = 1.033956264058147
$ws.Range("J2").Value = 1.025458205927374
$ws.Range("K2").Value = 1.034477217507755
$ws.Range("L2").Value = 1.024131211497256
$ws.Range("M2").Value = 1.035660871638634
$ws.Range("N2").Value = 1.012567783118517
# Row 3
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.021078495310792
$ws.Range("D3").Value = 1.032153822579851
$ws.Range("E3").Value = 1.021984120154224
$ws.Range("F3").Value = 1.033937492779029
$ws.Range("I3").Value = 1.034131683770662
$ws.Range("J3").Value = 1.025914406506355
$ws.Range("K3").Value = 1.034770753390998
$ws.Range("L3").Value = 1.02462861669534
$ws.Range("M3").Value = 1.036549649273495
$ws.Range("N3").Value = 1.012719492610585
# Row 4
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.021608504535737
$ws.Range("D4").Value = 1.032467076624351
$ws.Range("E4").Value = 1.02243181263032
$ws.Range("F4").Value = 1.034637389854166
$ws.Range("I4").Value = 1.034243908230484
$ws.Range("J4").Value = 1.026209401211766
$ws.Range("K4").Value = 1.034960008678433
$ws.Range("L4").Value = 1.024950581186568
$ws.Range("M4").Value = 1.037124805344475
$ws.Range("N4").Value = 1.012817563630214
# Row 5
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.021831425778745
$ws.Range("D5").Value = 1.032598776161908
$ws.Range("E5").Value = 1.022620208937197
$ws.Range("F5").Value = 1.03493178081192
$ws.Range("I5").Value = 1.034290779278061
$ws.Range("J5").Value = 1.026333368642177
$ws.Range("K5").Value = 1.035039407144898
$ws.Range("L5").Value = 1.025085960233548
$ws.Range("M5").Value = 1.037366613843545
$ws.Range("N5").Value = 1.012858769519752
# Row 6
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.021868861374251
$ws.Range("D6").Value = 1.032620889483885
$ws.Range("E6").Value = 1.022651852388957
$ws.Range("F6").Value = 1.034981219368441
$ws.Range("I6").Value = 1.034298631042899
$ws.Range("J6").Value = 1.026354180453955
$ws.Range("K6").Value = 1.03505272881888
$ws.Range("L6").Value = 1.025108692413288
$ws.Range("M6").Value = 1.037407215284688
$ws.Range("N6").Value = 1.012865686800209
# Row 7
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.021611482808039
$ws.Range("D7").Value = 1.032468836372652
$ws.Range("E7").Value = 1.022434329260647
$ws.Range("F7").Value = 1.034641322915604
$ws.Range("I7").Value = 1.034244535735358
$ws.Range("J7").Value = 1.026211057861907
$ws.Range("K7").Value = 1.034961070252306
$ws.Range("L7").Value = 1.024952390030578
$ws.Range("M7").Value = 1.03712803635194
$ws.Range("N7").Value = 1.012818114316841
# Row 8
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.020536580425482
$ws.Range("D8").Value = 1.031833347950771
$ws.Range("E8").Value = 1.021526704319336
$ws.Range("F8").Value = 1.033221922457708
$ws.Range("I8").Value = 1.034015813335729
$ws.Range("J8").Value = 1.0256124213307
$ws.Range("K8").Value = 1.03457655985563
$ws.Range("L8").Value = 1.024299288416359
$ws.Range("M8").Value = 1.035961225763057
$ws.Range("N8").Value = 1.012619073525717
# Row 9
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.018645942673253
$ws.Range("D9").Value = 1.030713912596411
$ws.Range("E9").Value = 1.019933439354176
$ws.Range("F9").Value = 1.030725756112145
$ws.Range("I9").Value = 1.033602976434558
$ws.Range("J9").Value = 1.02455608014577
$ws.Range("K9").Value = 1.033893831247065
$ws.Range("L9").Value = 1.023149334774158
$ws.Range("M9").Value = 1.033905635692423
$ws.Range("N9").Value = 1.012267626299466
# Row 10
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.017387938515461
$ws.Range("D10").Value = 1.029967981180431
$ws.Range("E10").Value = 1.018875450858292
$ws.Range("F10").Value = 1.029065054659025
$ws.Range("I10").Value = 1.033321208931085
$ws.Range("J10").Value = 1.023850927993524
$ws.Range("K10").Value = 1.033435270796361
$ws.Range("L10").Value = 1.022383370303658
$ws.Range("M10").Value = 1.032535615347312
$ws.Range("N10").Value = 1.012032872215588
# Row 11
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.016843800330387
$ws.Range("D11").Value = 1.029645092072365
$ws.Range("E11").Value = 1.0184183437693
$ws.Range("F11").Value = 1.028346770932678
$ws.Range("I11").Value = 1.033197658057104
$ws.Range("J11").Value = 1.023545382897388
$ws.Range("K11").Value = 1.033235915319168
$ws.Range("L11").Value = 1.022051873056894
$ws.Range("M11").Value = 1.03194248010559
$ws.Range("N11").Value = 1.01193111771959
# Row 12
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.016641772575129
$ws.Range("D12").Value = 1.029525174130285
$ws.Range("E12").Value = 1.018248707055057
$ws.Range("F12").Value = 1.028080090919585
$ws.Range("I12").Value = 1.033151534576688
$ws.Range("J12").Value = 1.023431859381655
$ws.Range("K12").Value = 1.03316174742907
$ws.Range("L12").Value = 1.021928766991058
$ws.Range("M12").Value = 1.031722177938802
$ws.Range("N12").Value = 1.011893306257583
# Row 13
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.016685104147182
$ws.Range("D13").Value = 1.029550896117388
$ws.Range("E13").Value = 1.01828508774678
$ws.Range("F13").Value = 1.028137289130873
$ws.Range("I13").Value = 1.033161438669584
$ws.Range("J13").Value = 1.023456211923253
$ws.Range("K13").Value = 1.033177662031951
$ws.Range("L13").Value = 1.021955172429884
$ws.Range("M13").Value = 1.031769432825297
$ws.Range("N13").Value = 1.011901417631363
# Row 14
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.016827098821582
$ws.Range("D14").Value = 1.029635179251483
$ws.Range("E14").Value = 1.018404318410682
$ws.Range("F14").Value = 1.02832472457563
$ws.Range("I14").Value = 1.033193850187664
$ws.Range("J14").Value = 1.023535999623687
$ws.Range("K14").Value = 1.033229786986614
$ws.Range("L14").Value = 1.022041696518213
$ws.Range("M14").Value = 1.031924269556072
$ws.Range("N14").Value = 1.011927992522515
# Row 15
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.016914598313446
$ws.Range("D15").Value = 1.029687111238199
$ws.Range("E15").Value = 1.018477800708611
$ws.Range("F15").Value = 1.028440225998844
$ws.Range("I15").Value = 1.033213789381474
$ws.Range("J15").Value = 1.023585155453195
$ws.Range("K15").Value = 1.033261887238586
$ws.Range("L15").Value = 1.022095010447059
$ws.Range("M15").Value = 1.032019671542989
$ws.Range("N15").Value = 1.011944364172388
# Row 16
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.017424063635979
$ws.Range("D16").Value = 1.029989412615172
$ws.Range("E16").Value = 1.018905808946487
$ws.Range("F16").Value = 1.029112741930577
$ws.Range("I16").Value = 1.033329376134893
$ws.Range("J16").Value = 1.023871201687279
$ws.Range("K16").Value = 1.033448484696423
$ws.Range("L16").Value = 1.022405374374481
$ws.Range("M16").Value = 1.032574981795633
$ws.Range("N16").Value = 1.012039623158428
# Row 17
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.017743795604107
$ws.Range("D17").Value = 1.030179067435474
$ws.Range("E17").Value = 1.019174558474987
$ws.Range("F17").Value = 1.029534811021968
$ws.Range("I17").Value = 1.033401467883278
$ws.Range("J17").Value = 1.024050575548761
$ws.Range("K17").Value = 1.03356532012501
$ws.Range("L17").Value = 1.022600103918358
$ws.Range("M17").Value = 1.032923338574036
$ws.Range("N17").Value = 1.012099348911982
# Row 18
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.017930346292414
$ws.Range("D18").Value = 1.030289699808542
$ws.Range("E18").Value = 1.019331412764037
$ws.Range("F18").Value = 1.029781074980558
$ws.Range("I18").Value = 1.033443368787779
$ws.Range("J18").Value = 1.024155180999368
$ws.Range("K18").Value = 1.033633391394936
$ws.Range("L18").Value = 1.022713702731971
$ws.Range("M18").Value = 1.0331265380735
$ws.Range("N18").Value = 1.012134175832426
# Row 19
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.017993964787981
$ws.Range("D19").Value = 1.030327424217143
$ws.Range("E19").Value = 1.019384912466271
$ws.Range("F19").Value = 1.029865057868785
$ws.Range("I19").Value = 1.033457630619243
$ws.Range("J19").Value = 1.02419084528248
$ws.Range("K19").Value = 1.033656588866708
$ws.Range("L19").Value = 1.02275243973204
$ws.Range("M19").Value = 1.033195825335606
$ws.Range("N19").Value = 1.012146049187464
# Row 20
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.017709485542025
$ws.Range("D20").Value = 1.030158718228481
$ws.Range("E20").Value = 1.019145714124626
$ws.Range("F20").Value = 1.029489518909281
$ws.Range("I20").Value = 1.033393748526627
$ws.Range("J20").Value = 1.0240313325221
$ws.Range("K20").Value = 1.033552792730684
$ws.Range("L20").Value = 1.022579209580032
$ws.Range("M20").Value = 1.032885962264046
$ws.Range("N20").Value = 1.012092941945941
# Row 21
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.016785282450386
$ws.Range("D21").Value = 1.029610359473038
$ws.Range("E21").Value = 1.018369203709424
$ws.Range("F21").Value = 1.02826952612033
$ws.Range("I21").Value = 1.033184312179342
$ws.Range("J21").Value = 1.023512504967284
$ws.Range("K21").Value = 1.033214440745427
$ws.Range("L21").Value = 1.022016216582366
$ws.Range("M21").Value = 1.031878673586957
$ws.Range("N21").Value = 1.011920167299077
# Row 22
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.016204717171612
$ws.Range("D22").Value = 1.029265686403985
$ws.Range("E22").Value = 1.017881868031252
$ws.Range("F22").Value = 1.027503177590158
$ws.Range("I22").Value = 1.033051294036037
$ws.Range("J22").Value = 1.023186121748908
$ws.Range("K22").Value = 1.033001021032264
$ws.Range("L22").Value = 1.021662396007836
$ws.Range("M22").Value = 1.031245436478141
$ws.Range("N22").Value = 1.011811448619588
# Row 23
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.016512436264021
$ws.Range("D23").Value = 1.0294483938644
$ws.Range("E23").Value = 1.018140129250303
$ws.Range("F23").Value = 1.027909365969772
$ws.Range("I23").Value = 1.033121935975828
$ws.Range("J23").Value = 1.023359159985872
$ws.Range("K23").Value = 1.033114223344054
$ws.Range("L23").Value = 1.021849947863042
$ws.Range("M23").Value = 1.03158111911496
$ws.Range("N23").Value = 1.011869090702142
# Row 24
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.017724988607972
$ws.Range("D24").Value = 1.030167913128764
$ws.Range("E24").Value = 1.019158747344721
$ws.Range("F24").Value = 1.029509984223796
$ws.Range("I24").Value = 1.033397237032063
$ws.Range("J24").Value = 1.02404002768034
$ws.Range("K24").Value = 1.033558453558567
$ws.Range("L24").Value = 1.022588650781703
$ws.Range("M24").Value = 1.032902850982242
$ws.Range("N24").Value = 1.012095837009364
# Row 25
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.019134297203413
$ws.Range("D25").Value = 1.031003257935918
$ws.Range("E25").Value = 1.020344605376579
$ws.Range("F25").Value = 1.031370477584109
$ws.Range("I25").Value = 1.033710861211819
$ws.Range("J25").Value = 1.024829337160908
$ws.Range("K25").Value = 1.034070938716807
$ws.Range("L25").Value = 1.023446512127201
$ws.Range("M25").Value = 1.034436992977851
$ws.Range("N25").Value = 1.012358565854002
